$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: remove the standalone "Meta description: ..." paragraph
# that used to sit right under the page title.
# ------------------------------------------------------------------
$metaFindRange = $d.Content
$metaFound = $metaFindRange.Find.Execute(
    "Meta description: Read our expert review of Battleship Direct Hit and play for free. Impressive graphics, Megaways, and free spins make for an exciting gaming experience.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($metaFound) {
    $metaPara = $metaFindRange.Paragraphs(1)
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# Change 2 & 3: at the bottom of the document, the paragraph that used
# to hold the AI image-generation prompt now instead holds the "meta
# description" copy, and a new bold heading-like paragraph ("Play
# Battleship Direct Hit Free - Exciting Online Slot Game") is inserted
# right before it.
# ------------------------------------------------------------------
$promptFindRange = $d.Content
$promptFound = $promptFindRange.Find.Execute(
    'Please create an image featuring a happy Maya warrior wearing glasses for the game "Battleship Direct Hit". The image should be in a cartoon style and should capture the essence of the game''s naval battle theme in a fun and engaging way. It should also feature the game''s title prominently. Be creative and use vibrant colors and dynamic imagery to attract players to this exciting slot game.',
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($promptFound) {
    $promptPara = $promptFindRange.Paragraphs(1)

    # --- insert a new, empty paragraph right before it ---
    $promptPara.Range.InsertParagraphBefore()
    $count = $d.Paragraphs.Count
    $newPara = $d.Paragraphs($count - 1)

    # The freshly inserted paragraph mark inherits formatting (e.g.
    # italics) from the paragraph that follows it. Neutralize that by
    # briefly typing a placeholder character, clearing its formatting,
    # then removing it again - leaving a clean, unformatted paragraph.
    $placeholderIns = $d.Range($newPara.Range.Start, $newPara.Range.Start)
    $placeholderIns.InsertBefore("X")
    $placeholderRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
    $placeholderRange.Font.Bold = 0
    $placeholderRange.Font.Italic = 0
    $placeholderRange2 = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
    $placeholderRange2.Delete()

    # Now fill the clean paragraph with an empty leading run followed
    # by a bold run, matching the structure used elsewhere in the doc.
    $newPara2 = $d.Paragraphs($count - 1)
    $newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Battleship Direct Hit Free - Exciting Online Slot Game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newParaTarget = $d.Range($newPara2.Range.Start, $newPara2.Range.Start)
    $newParaTarget.InsertXML($newParaXml)

    # --- swap out the old AI-image-prompt text for the new meta copy,
    #     keeping the paragraph's existing (italic) run formatting ---
    $d.Content.Find.Execute(
        'Please create an image featuring a happy Maya warrior wearing glasses for the game "Battleship Direct Hit". The image should be in a cartoon style and should capture the essence of the game''s naval battle theme in a fun and engaging way. It should also feature the game''s title prominently. Be creative and use vibrant colors and dynamic imagery to attract players to this exciting slot game.',
        $true, $false, $false, $false, $false, $true, 1, $false,
        "Read our expert review of Battleship Direct Hit and play for free. Impressive graphics, Megaways, and free spins make for an exciting gaming experience.",
        2)
}
